$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "93×24=2232" "67×68=4556"
Replace-Text "33×19=627" "81×69=5589"
Replace-Text "82×89=7298" "68×38=2584"
Replace-Text "56×65=3640" "48×41=1968"
Replace-Text "62×30=1860" "93×92=8556"
Replace-Text "38×73=2774" "72×36=2592"
Replace-Text "52×11=572" "15×90=1350"
Replace-Text "80×27=2160" "54×92=4968"
Replace-Text "18×73=1314" "75×72=5400"
Replace-Text "64×92=5888" "32×13=416"
Replace-Text "67×66=4422" "54×77=4158"
Replace-Text "94×53=4982" "86×45=3870"
Replace-Text "16×90=1440" "25×55=1375"
Replace-Text "96×59=5664" "40×17=680"
Replace-Text "98×34=3332" "24×69=1656"
Replace-Text "71×76=5396" "68×56=3808"
Replace-Text "74×51=3774" "75×50=3750"
Replace-Text "76×69=5244" "51×57=2907"
Replace-Text "33×56=1848" "74×41=3034"
Replace-Text "55×18=990" "77×51=3927"
Replace-Text "74×40=2960" "57×38=2166"
Replace-Text "27×89=2403" "46×63=2898"
Replace-Text "58×73=4234" "61×83=5063"
Replace-Text "52×44=2288" "30×89=2670"
Replace-Text "29×38=1102" "18×24=432"
